$d = $word.ActiveDocument

# Locate the three paragraphs to remove by their content:
#   1) the blank paragraph right after the "Bibliografia" body text
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "(c) 2020 ... Creative Commons Attribution" footer line
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startPara = $target.Previous()   # the blank paragraph just before it
    $endPara = $target.Next()         # the copyright paragraph just after it

    $start = $startPara.Range.Start
    $end = $endPara.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
